# Rename the workbook's first/active sheet from "Sheet1" to "ProgrammingBooks".
# This supports reading the sheet with the Poiji library into a POJO
# (e.g. a "ProgrammingBooks" mapped class) as noted in the commit message.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Name = "ProgrammingBooks"
